$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 49, shifting rows 49-50 down to 50-51,
# and copying the formatting of the existing row 49.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new data entry.
$ws.Range("A49").Value = 6
$ws.Range("B49").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C49").Value = "Metropolitana"
$ws.Range("D49").Value = 44448
$ws.Range("D49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E49").Value = 13
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100108
$ws.Range("H49").Value = "Tropicales y subtropicales"
$ws.Range("I49").Value = 100108007
$ws.Range("J49").Value = "Coco"
$ws.Range("K49").Value = "Sin especificar"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 150
$ws.Range("N49").Value = 20000
$ws.Range("O49").Value = 20000
$ws.Range("P49").Value = 20000
$ws.Range("Q49").Value = "$/malla 20 unidades"
$ws.Range("R49").Value = "Perú"
$ws.Range("S49").Value = 1000
$ws.Range("T49").Value = 20
